$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

$ws.Range("A7").Value = "Elmar Qara"
$ws.Range("B7").Value = "elmarqarayev69@gmail.com"
$ws.Range("C7").Value = 1794
$ws.Range("D7").Value = "Pending"
